# Update "想去人数" (number of people interested) figures for a handful of
# events. The same underlying records appear on two sheets: "展览" (only the
# exhibition rows) and "全部类型" (all rows merged together, so the same
# events land on different row numbers).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 281
$ws1.Range("F3").Value = 1161
$ws1.Range("F4").Value = 16589
$ws1.Range("F6").Value = 1626
$ws1.Range("F8").Value = 358
$ws1.Range("F9").Value = 206
$ws1.Range("F11").Value = 11544
$ws1.Range("F13").Value = 1221
$ws1.Range("F14").Value = 4568
$ws1.Range("F15").Value = 401
$ws1.Range("F19").Value = 330
$ws1.Range("F20").Value = 147

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 281
$ws4.Range("F4").Value = 1161
$ws4.Range("F5").Value = 16589
$ws4.Range("F7").Value = 1626
$ws4.Range("F9").Value = 358
$ws4.Range("F10").Value = 206
$ws4.Range("F14").Value = 11544
$ws4.Range("F16").Value = 1221
$ws4.Range("F17").Value = 4568
$ws4.Range("F18").Value = 401
$ws4.Range("F22").Value = 330
$ws4.Range("F23").Value = 147
